$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix inconsistent capitalization of "RData" -> "Rdata " (matches the rest of the column)
$ws.Range("A3").Value = "Rdata "

# These four rows now share the rerun pipeline name instead of individual timestamps
$ws.Range("B14").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("B15").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("B16").Value = "Pipe_SCTv2_corrected_13-06"
$ws.Range("B17").Value = "Pipe_SCTv2_corrected_13-06"

# New Kriegstein-to-SingleR / Rdata rows logging the latest SCTv2-corrected runs
$ws.Range("A19").Value = "Rdata "
$ws.Range("B19").Value = "SingleR_RData_2022-06-13 14-46-51"
$ws.Range("C19").Value = "Kriegstein to SingleR"
$ws.Range("D19").Value = "SCTv2 corrected new selection"
$ws.Range("F19").Value = "rerun SCTv2 corrected pipeline (integrated)"
$ws.Range("G19").Value = "DEG, pseudotime"

$ws.Range("A20").Value = "Rdata "
$ws.Range("B20").Value = "SingleR_RData_2022-06-13 14-50-40"
$ws.Range("C20").Value = "Kriegstein to SingleR"
$ws.Range("D20").Value = "SCTv2 corrected new post selection"
$ws.Range("F20").Value = "rerun SCTv2 corrected pipeline (integrated)"
$ws.Range("G20").Value = "DEG, pseudotime"

$ws.Range("A21").Value = "Rdata "
$ws.Range("B21").Value = "SingleR_RData_2022-06-13 14-55-08"
$ws.Range("C21").Value = "Kriegstein to SingleR"
$ws.Range("D21").Value = "SCTv2 corrected old selection"
$ws.Range("F21").Value = "rerun SCTv2 corrected pipeline (individual + integrated)"
$ws.Range("G21").Value = "DEG, pseudotime"

$ws.Range("A22").Value = "Rdata "
$ws.Range("B22").Value = "SingleR_RData_2022-06-13 14-56-28"
$ws.Range("C22").Value = "Kriegstein to SingleR"
$ws.Range("D22").Value = "SCTv2 corrected old post selection"
$ws.Range("F22").Value = "rerun SCTv2 corrected pipeline (integrated)"
$ws.Range("G22").Value = "DEG, pseudotime"

# Scroll/selection moved to B25 (and the frozen top-left cell is cleared)
$ws.Range("B25").Select()
